$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The weekly report table (Table1) currently ends at row 61. Add a new row
# for 21/7/2025 the same way the author would: grow the table by one row
# (this keeps the table ref / autoFilter / dimension in sync), then copy
# the formatting of the previous week's row onto it and fill in this
# week's numbers.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.ListRows.Add() | Out-Null

$newRowIndex = 62
$prevRowIndex = $newRowIndex - 1

$srcRange = $ws.Range("D" + $prevRowIndex + ":J" + $prevRowIndex)
$dstRange = $ws.Range("D" + $newRowIndex + ":J" + $newRowIndex)
$srcRange.Copy($dstRange)
$ws.Rows.Item($newRowIndex).RowHeight = $ws.Rows.Item($prevRowIndex).RowHeight

$ws.Cells.Item($newRowIndex, 4).Value = "21/7/2025"   # Fecha
$ws.Cells.Item($newRowIndex, 5).Value = 406            # Imagenes sin etiquetar
$ws.Cells.Item($newRowIndex, 6).Value = 924            # Imagenes etiquetadas sin revisar
$ws.Cells.Item($newRowIndex, 7).Value = 0              # Imagenes rechazadas
$ws.Cells.Item($newRowIndex, 8).Value = 0              # Imagenes etiquetadas y revisadas, faltando de subir
$ws.Cells.Item($newRowIndex, 9).Value = 1012           # Imagenes etiquetadas, revisadas y subidas
$ws.Cells.Item($newRowIndex, 10).Value = "N/A"         # Notas

# Scroll the sheet down and select the cell just below the new row, mirroring
# the view state saved by the author after adding this entry.
$win = $excel.ActiveWindow
$win.ScrollRow = 52
$win.ScrollColumn = 2
$ws.Range("D" + ($newRowIndex + 1)).Select() | Out-Null
